$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Tool Log: record Kevin Johnson checking out / returning a Clamp
# ---------------------------------------------------------------------
$toolLog = $wb.Worksheets.Item("Tool Log")
$toolLog.Range("A3").Value = "Kevin Johnson has checked out 4 of Clamp"
$toolLog.Range("B3").Value = 42131.733229166668
$toolLog.Range("B1").Copy()
$toolLog.Range("B3").PasteSpecial(-4122)

$toolLog.Range("A4").Value = "Kevin Johnson has returned 2 of Clamp"
$toolLog.Range("B4").Value = 42131.73400462963
$toolLog.Range("B1").Copy()
$toolLog.Range("B4").PasteSpecial(-4122)

$toolLog.Range("A4:B4").Select()

# ---------------------------------------------------------------------
# 2. Purchase Log: record Kevin Johnson purchasing / returning a Limit Switch
# ---------------------------------------------------------------------
$purchaseLog = $wb.Worksheets.Item("Purchase Log")
$purchaseLog.Range("A3").Value = "Kevin Johnson has purchased 5 of Limit Switch"
$purchaseLog.Range("B3").Value = 42131.735960648148
$purchaseLog.Range("B1").Copy()
$purchaseLog.Range("B3").PasteSpecial(-4122)

$purchaseLog.Range("A4").Value = "Kevin Johnson has returned 3 of Limit Switch"
$purchaseLog.Range("B4").Value = 42131.736712962964
$purchaseLog.Range("B1").Copy()
$purchaseLog.Range("B4").PasteSpecial(-4122)

$purchaseLog.Range("A4:B4").Select()

# ---------------------------------------------------------------------
# 3. Student: add new student Julian Patrick at Encinal
# ---------------------------------------------------------------------
$student = $wb.Worksheets.Item("Student")
$student.Range("A316").Value = "Julian Patrick"
$student.Range("B316").Value = "Encinal"
$student.Range("A316:B316").Select()

# ---------------------------------------------------------------------
# 4. Tool: add new tool Claw Hammer, update Clamp quantity
# ---------------------------------------------------------------------
$tool = $wb.Worksheets.Item("Tool")
$tool.Range("B17").Value = 5
$tool.Range("A36").Value = "Claw Hammer"
$tool.Range("B36").Value = 4
$tool.Range("A36:B36").Select()

# ---------------------------------------------------------------------
# 5. School: update Bishop O'Dowd tool count
# ---------------------------------------------------------------------
$school = $wb.Worksheets.Item("School")
$school.Range("B7").Value = 238
$school.Range("B7").Select()

# ---------------------------------------------------------------------
# 6. Purchase: clear a stale RANDBETWEEN value, add new purchase item
# ---------------------------------------------------------------------
$purchase = $wb.Worksheets.Item("Purchase")
$purchase.Range("C45").Value = 78
$purchase.Range("A57").Value = "VEX Controller"
$purchase.Range("B57").Value = 30
$purchase.Range("C57").Value = 25
$purchase.Range("A57:C57").Select()

Write-Output "done"
